$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 08:52"

# Row 25 - Israel: refreshed stats
$ws.Range("B25").Value = 14326
$ws.Range("C25").Value = 384
$ws.Range("D25").Value = 4961
$ws.Range("E25").Value = 9178
$ws.Range("F25").Value = 148
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 187

# Row 40 - Chequia: refreshed stats
$ws.Range("B40").Value = 7041
$ws.Range("C40").Value = 8
$ws.Range("E40").Value = 5087

# Rows 45-46 - Ucrania overtakes Catar in the ranking
$ws.Range("A45").Value = "Ucrania"
$ws.Range("B45").Value = 6592
$ws.Range("C45").Value = 467
$ws.Range("D45").Value = 424
$ws.Range("E45").Value = 5994
$ws.Range("F45").Value = 45
$ws.Range("G45").Value = 13
$ws.Range("H45").Value = 174

$ws.Range("A46").Value = "Catar"
$ws.Range("B46").Value = 6533
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 614
$ws.Range("E46").Value = 5910
$ws.Range("F46").Value = 37
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 9

# Rows 69-71 - Oman overtakes Irak and Estonia in the ranking
$ws.Range("A69").Value = "Oman"
$ws.Range("B69").Value = 1614
$ws.Range("C69").Value = 106
$ws.Range("D69").Value = 238
$ws.Range("E69").Value = 1368
$ws.Range("F69").Value = 3
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 8

$ws.Range("A70").Value = "Irak"
$ws.Range("B70").Value = 1602
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 1096
$ws.Range("E70").Value = 423
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 83

$ws.Range("A71").Value = "Estonia"
$ws.Range("B71").Value = 1552
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 169
$ws.Range("E71").Value = 1340
$ws.Range("F71").Value = 9
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 43

# Row 122 - Vietnam: refreshed stats
$ws.Range("D122").Value = 222
$ws.Range("E122").Value = 46
